$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step date moves forward by one day
$ws.Range("A1").Value = 45309

# Updated price list (step 1 and 2 price increase)
$ws.Range("D19").Value = 2830.888
$ws.Range("D20").Value = 2997.795
$ws.Range("D21").Value = 3365.72
$ws.Range("D22").Value = 3733.619
$ws.Range("D23").Value = 4176.48
$ws.Range("D24").Value = 4510.331
$ws.Range("D25").Value = 4994.066
$ws.Range("D26").Value = 5300.658
$ws.Range("D27").Value = 5859.339
$ws.Range("D28").Value = 6581.534
$ws.Range("D29").Value = 7290.115
$ws.Range("D30").Value = 8243.959000000001
$ws.Range("D31").Value = 9674.727999999999

$ws.Range("D38").Value = 9960.880999999999
$ws.Range("D39").Value = 11745.938
$ws.Range("D40").Value = 13694.495
$ws.Range("D41").Value = 17312.306
$ws.Range("D42").Value = 22258.673
$ws.Range("D43").Value = 27402.637
$ws.Range("D44").Value = 31000.001
$ws.Range("D45").Value = 35632.961
$ws.Range("D46").Value = 15656.693

$ws.Range("D53").Value = 4898.677
$ws.Range("D54").Value = 5886.587
$ws.Range("D55").Value = 6315.816
$ws.Range("D56").Value = 7167.479
$ws.Range("D57").Value = 7796.208
$ws.Range("D58").Value = 8414.282999999999
$ws.Range("D59").Value = 9122.858
$ws.Range("D60").Value = 9409.008
$ws.Range("D61").Value = 10274.275
$ws.Range("D62").Value = 12181.985
$ws.Range("D63").Value = 13149.458
$ws.Range("D64").Value = 16215.381
